$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns stay text-formatted so values are not
# auto-converted to numbers/dates by Excel when assigned.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 20 and row 21 swap places (Chainlink <-> Uniswap) with updated data
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "10.18"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "15.25"
$ws.Range("E21").Value = "  -4.80%  "

$ws.Range("D2").Value = "66.956.56"
$ws.Range("E2").Value = "  -1.97%  "

$ws.Range("D3").Value = "3.483.58"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "601.53"
$ws.Range("E5").Value = "  -2.89%  "

$ws.Range("D6").Value = "148.04"
$ws.Range("E6").Value = "  -4.59%  "

$ws.Range("D7").Value = "3.480.56"
$ws.Range("E7").Value = "  -2.27%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("E10").Value = "  -2.53%  "

$ws.Range("D11").Value = "7.66"
$ws.Range("E11").Value = "  +3.58%  "

$ws.Range("D12").Value = "0.423"
$ws.Range("E12").Value = "  -3.44%  "

$ws.Range("E13").Value = "  -3.65%  "

$ws.Range("D14").Value = "4.068.04"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "31.19"
$ws.Range("E15").Value = "  -6.28%  "

$ws.Range("D16").Value = "3.474.47"
$ws.Range("E16").Value = "  -2.60%  "

$ws.Range("D17").Value = "66.936.25"
$ws.Range("E17").Value = "  -2.08%  "

$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  -5.02%  "

$ws.Range("D22").Value = "434.60"
$ws.Range("E22").Value = "  -4.32%  "

$ws.Range("D23").Value = "0.606"
$ws.Range("E23").Value = "  -5.78%  "

$ws.Range("D24").Value = "79.12"
$ws.Range("E24").Value = "  +0.83%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").Value = "3.616.77"
$ws.Range("E26").Value = "  -2.46%  "

$ws.Range("E27").Value = "  -9.50%  "

$ws.Range("D28").Value = "9.82"
$ws.Range("E28").Value = "  -6.99%  "

$ws.Range("D29").Value = "8.36"
$ws.Range("E29").Value = "  -9.19%  "

$ws.Range("D30").Value = "2.48"
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").Value = "1.58"
$ws.Range("E31").Value = "  -7.38%  "

$ws.Range("E32").Value = "  -2.45%  "

$ws.Range("E33").Value = "  -0.10%  "

$ws.Range("D34").Value = "25.34"
$ws.Range("E34").Value = "  -3.13%  "

$ws.Range("D35").Value = "3.469.66"
$ws.Range("E35").Value = "  -2.49%  "

$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  -6.05%  "

$ws.Range("D37").Value = "5.91"
$ws.Range("E37").Value = "  -7.36%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").Value = "7.89"
$ws.Range("E39").Value = "  -4.15%  "

$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "173.54"
$ws.Range("E41").Value = "  -4.22%  "

$ws.Range("D42").Value = "0.0884"
$ws.Range("E42").Value = "  -3.63%  "

$ws.Range("E43").Value = "  -12.94%  "

$ws.Range("D44").Value = "5.40"
$ws.Range("E44").Value = "  -3.68%  "

$ws.Range("D45").Value = "0.897"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "46.45"
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").Value = "28.84"
$ws.Range("E47").Value = "  -7.27%  "

$ws.Range("E48").Value = "  -7.19%  "

$ws.Range("D49").Value = "7.46"
$ws.Range("E49").Value = "  -4.26%  "

$ws.Range("D50").Value = "2.42"
$ws.Range("E50").Value = "  -9.42%  "

$ws.Range("E51").Value = "  -4.79%  "
